# 17.1.1.xlsx update:
#  - add a new "2020" data column (Q) to the revenue table
#  - shrink the (now less-tall) header row
#  - leave the selection where the author left it (N13, just below the table)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new "2020" column ------------------------------------------------
# Year header, matching the style already used for the other year headers
# in row 4 (bold Times New Roman, right/center, top+bottom border).
$ws.Range("Q4").Value = 2020
$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)   # xlPasteFormats

# "Revenues, total" row - bold summary row, needs its own (new) style:
# bold, 9pt Times New Roman, number format 0.0, right/center aligned.
$ws.Range("Q5").Value = 25.3
$ws.Range("Q5").Font.Name = "Times New Roman"
$ws.Range("Q5").Font.Family = 1
$ws.Range("Q5").Font.Size = 9
$ws.Range("Q5").Font.ColorIndex = -4105
$ws.Range("Q5").NumberFormat = "0.0"
$ws.Range("Q5").HorizontalAlignment = -4152   # xlRight
$ws.Range("Q5").VerticalAlignment = -4108     # xlCenter

# Remaining data rows - copy the formatting already used a couple of
# columns over (N) so the new cells render exactly like their neighbours.
$ws.Range("Q6").Value = 17.8
$ws.Range("N6").Copy()
$ws.Range("Q6").PasteSpecial(-4122)

$ws.Range("Q7").Value = "-"
$ws.Range("N7").Copy()
$ws.Range("Q7").PasteSpecial(-4122)

$ws.Range("Q8").Value = 2
$ws.Range("N8").Copy()
$ws.Range("Q8").PasteSpecial(-4122)

$ws.Range("Q9").Value = 5.5
$ws.Range("P9").Copy()
$ws.Range("Q9").PasteSpecial(-4122)

$ws.Range("Q10").Value = 0
$ws.Range("P10").Copy()
$ws.Range("Q10").PasteSpecial(-4122)

# --- header row height --------------------------------------------------
$ws.Rows.Item(1).RowHeight = 38.25

# --- restore author's last selection ------------------------------------
$ws.Range("N13").Select()
